$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove H2 cell (formula + value)
$ws.Range("H2").ClearContents()

# G3 becomes a dynamic array formula spilling into G3:J3
$ws.Range("G3:J3").FormulaArray = "=TEXTSPLIT(G2,"","","","",1)"
